# Phase 0 section: highlight everything green, and merge the "Done when"
# run that was split across two <w:r> elements into a single run.
#
# wdBrightGreen = 4  ->  OOXML <w:highlight w:val="green"/>
$wdBrightGreen = 4

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the two runs of the "Done when: ... UI" / " and blocks ..."
#    paragraph (numId=2) into a single run BEFORE applying highlighting,
#    so the highlight ends up on one merged run as in the target XML.
# ---------------------------------------------------------------------
$mergeText = "Done when: refresh keeps user logged in; logout clears UI and blocks purchase/offer."
$doneParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd("`r") -eq $mergeText) {
        $doneParaIndex = $i
        break
    }
}

if ($doneParaIndex -gt 0) {
    $donePara = $d.Paragraphs.Item($doneParaIndex)
    $pRange = $donePara.Range
    $textRange = $d.Range($pRange.Start, $pRange.End - 1)
    $textRange.Delete()
    $insPoint = $d.Range($pRange.Start, $pRange.Start)
    $insPoint.InsertBefore($mergeText)
}

# ---------------------------------------------------------------------
# 2) Highlight the Phase 0 paragraphs (heading through the end of the
#    "Guest flow enforcement" sub-section) green. Using Font on the full
#    paragraph Range (which includes the paragraph mark) sets both the
#    run rPr *and* the paragraph mark's pPr/rPr highlight, matching the
#    target markup.
# ---------------------------------------------------------------------
$targets = @(
    "Phase 0 — Stabilize Auth & Sessions (foundational)",
    "Why now: Everything else sits on this.",
    "Session & Logout",
    "Implement: supabase.auth.getSession() on app load; store user in state; add logout() button.",
    "Guard: purchase/offer actions open login if no session.",
    $mergeText,
    "Guest flow enforcement",
    "Gate " + [char]0x201C + "Buy" + [char]0x201D + " & " + [char]0x201C + "Offer price" + [char]0x201D + " " + [char]0x2192 + " if guest, open login modal/page (keep return" + [char]0x001E + "to).",
    "Done when: guests can browse but cannot complete protected actions."
)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $ptext = $para.Range.Text.TrimEnd("`r")
    if ($targets -contains $ptext) {
        $para.Range.Font.HighlightColorIndex = $wdBrightGreen
    }
}

Write-Output "Phase 0 highlighting applied."
